$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A10 used to hold the shared string "foo" - replace it with the numeric value 0
$ws.Range("A10").Value = 0

# D10 changes from 10 to 0
$ws.Range("D10").Value = 0

# Update the active selection shown when the sheet is opened
$ws.Range("E10").Select()
